$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "Save" column (match the header style used by G1, e.g.)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Save values per row (2-30), matching the commit's H column data
$saveValues = @{
    2 = 1; 3 = 0; 4 = 0; 5 = 0; 6 = 1; 7 = 0; 8 = 0; 9 = 1; 10 = 1;
    11 = 0; 12 = 1; 13 = 0; 14 = 0; 15 = 0; 16 = 0; 17 = 1; 18 = 0; 19 = 1;
    20 = 0; 21 = 0; 22 = 0; 23 = 0; 24 = 0; 25 = 1; 26 = 0; 27 = 0; 28 = 0;
    29 = 0; 30 = 1
}

foreach ($row in 2..30) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
